# Update scripts.xlsx: add two new script entries to the inventory sheet,
# size the columns to fit their content, and leave the active selection
# on the row just below the newly entered data (matches typical Excel
# behaviour after typing a block of rows and pressing Enter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
$ws.Range("A2").Value = "script000.py"
$ws.Range("B2").Value = "Arquivo contendo funções gerais úteis nos demais scripts"

$ws.Range("A3").Value = "script001.py"
$ws.Range("B3").Value = "Carrega os dados de Itacoatiara e salva esses dados em arquivos separados"

# --- Column widths ---------------------------------------------------------
# Column A was auto-fit to the longest entry ("NOME DO ARQUIVO" / "script001.py"),
# columns B and C were sized by hand to comfortably fit the descriptions.
$ws.Columns.Item(1).ColumnWidth = 17.333333333333332
$ws.Columns.Item(2).ColumnWidth = 96.16666666666667
$ws.Columns.Item(3).ColumnWidth = 62

# --- Selection -------------------------------------------------------------
# After entering the two rows, the cursor rests one row below the data.
$ws.Range("A4").Select() | Out-Null
